$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.076880890253462
$ws.Range("D2").Value = 1.075313900309295
$ws.Range("E2").Value = 1.0806092323852
$ws.Range("F2").Value = 1.084210506821714
$ws.Range("I2").Value = 1.049269138000093
$ws.Range("J2").Value = 1.081778133705163
$ws.Range("K2").Value = 1.078000798051144
$ws.Range("L2").Value = 1.083282219649911
$ws.Range("M2").Value = 1.086874117584125
$ws.Range("N2").Value = 1.083314382736711
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.07907451745983
$ws.Range("D3").Value = 1.077050357833876
$ws.Range("E3").Value = 1.082757062749616
$ws.Range("F3").Value = 1.08608631845749
$ws.Range("I3").Value = 1.049872422094993
$ws.Range("J3").Value = 1.083626000076901
$ws.Range("K3").Value = 1.07955229591135
$ws.Range("L3").Value = 1.08524511918176
$ws.Range("M3").Value = 1.088566348822414
$ws.Range("N3").Value = 1.08516487329065
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.080488214116116
$ws.Range("D4").Value = 1.078168599514492
$ws.Range("E4").Value = 1.084141497781082
$ws.Range("F4").Value = 1.087294869564981
$ws.Range("I4").Value = 1.050258725365911
$ws.Range("J4").Value = 1.084815810605727
$ws.Range("K4").Value = 1.080550401978188
$ws.Range("L4").Value = 1.086509518702397
$ws.Range("M4").Value = 1.08965567998693
$ws.Range("N4").Value = 1.0863563734869
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.081081193339376
$ws.Range("D5").Value = 1.078637449232347
$ws.Range("E5").Value = 1.084722262813006
$ws.Range("F5").Value = 1.087801719831912
$ws.Range("I5").Value = 1.050420165055844
$ws.Range("J5").Value = 1.085314626112974
$ws.Range("K5").Value = 1.080968636609881
$ws.Range("L5").Value = 1.087039729724541
$ws.Range("M5").Value = 1.090112305549905
$ws.Range("N5").Value = 1.086855897369382
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.081180679578642
$ws.Range("D6").Value = 1.078716098043178
$ws.Range("E6").Value = 1.08481970331277
$ws.Range("F6").Value = 1.087886751164579
$ws.Range("I6").Value = 1.050447215361052
$ws.Range("J6").Value = 1.085398299251852
$ws.Range("K6").Value = 1.081038780432517
$ws.Range("L6").Value = 1.087128676590253
$ws.Range("M6").Value = 1.090188897756504
$ws.Range("N6").Value = 1.086939689333714
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.080496142751052
$ws.Range("D7").Value = 1.078174869221447
$ws.Range("E7").Value = 1.084149262868198
$ws.Range("F7").Value = 1.087301646899784
$ws.Range("I7").Value = 1.050260886297961
$ws.Range("J7").Value = 1.084822481191752
$ws.Range("K7").Value = 1.080555995798161
$ws.Range("L7").Value = 1.086516608644672
$ws.Range("M7").Value = 1.08966178662986
$ws.Range("N7").Value = 1.086363053545922
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.077623442102479
$ws.Range("D8").Value = 1.07590186955564
$ws.Range("E8").Value = 1.081336228947293
$ws.Range("F8").Value = 1.084845543733783
$ws.Range("I8").Value = 1.049473867916275
$ws.Range("J8").Value = 1.082403863462941
$ws.Range("K8").Value = 1.078526352598675
$ws.Range("L8").Value = 1.083946794087095
$ws.Range("M8").Value = 1.087447201053424
$ws.Range("N8").Value = 1.083941001102517
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.07251605898522
$ws.Range("D9").Value = 1.071854364604896
$ws.Range("E9").Value = 1.076336910430344
$ws.Range("F9").Value = 1.080476379344777
$ws.Range("I9").Value = 1.048055468413147
$ws.Range("J9").Value = 1.078095652235443
$ws.Range("K9").Value = 1.074904252025384
$ws.Range("L9").Value = 1.079373284950467
$ws.Range("M9").Value = 1.083500383649375
$ws.Range("N9").Value = 1.079626671720917
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.069078613903569
$ws.Range("D10").Value = 1.069126083403361
$ws.Range("E10").Value = 1.072973580286818
$ws.Range("F10").Value = 1.07753423486328
$ws.Range("I10").Value = 1.047087986434537
$ws.Range("J10").Value = 1.075190628389317
$ws.Range("K10").Value = 1.072457356700653
$ws.Range("L10").Value = 1.076292077384321
$ws.Range("M10").Value = 1.080837717193419
$ws.Range("N10").Value = 1.076717522407719
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.067581986235106
$ws.Range("D11").Value = 1.067937249623815
$ws.Range("E11").Value = 1.071509560367606
$ws.Range("F11").Value = 1.076252912258799
$ws.Range("I11").Value = 1.046663718209462
$ws.Range("J11").Value = 1.073924526639146
$ws.Range("K11").Value = 1.07138985886229
$ws.Range("L11").Value = 1.074949830232743
$ws.Range("M11").Value = 1.07967693698378
$ws.Range("N11").Value = 1.075449622647839
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.067024801334658
$ws.Range("D12").Value = 1.067494510646107
$ws.Range("E12").Value = 1.070964566939652
$ws.Range("F12").Value = 1.075775834389993
$ws.Range("I12").Value = 1.046505310626181
$ws.Range("J12").Value = 1.07345297271668
$ws.Range("K12").Value = 1.070992115039628
$ws.Range("L12").Value = 1.074450012533083
$ws.Range("M12").Value = 1.079244564544883
$ws.Range("N12").Value = 1.074977399064715
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.067144377542221
$ws.Range("D13").Value = 1.067589532345102
$ws.Range("E13").Value = 1.071081524403947
$ws.Range("F13").Value = 1.075878221221361
$ws.Range("I13").Value = 1.046539326701543
$ws.Range("J13").Value = 1.073554180541509
$ws.Range("K13").Value = 1.071077488466563
$ws.Range("L13").Value = 1.074557282152927
$ws.Range("M13").Value = 1.079337365021172
$ws.Range("N13").Value = 1.075078750616263
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.067535955314745
$ws.Range("D14").Value = 1.067900676386953
$ws.Range("E14").Value = 1.071464535564313
$ws.Range("F14").Value = 1.076213500251006
$ws.Range("I14").Value = 1.046650640899946
$ws.Range("J14").Value = 1.073885573899186
$ws.Range("K14").Value = 1.071357006447669
$ws.Range("L14").Value = 1.074908540800695
$ws.Range("M14").Value = 1.079641221745959
$ws.Range("N14").Value = 1.07541061459052
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.067777049706251
$ws.Range("D15").Value = 1.06809222875936
$ws.Range("E15").Value = 1.071700362380317
$ws.Range("F15").Value = 1.07641992499689
$ws.Range("I15").Value = 1.046719116804069
$ws.Range("J15").Value = 1.074089587310153
$ws.Range("K15").Value = 1.07152906314787
$ws.Range("L15").Value = 1.075124796478973
$ws.Range("M15").Value = 1.079828277031565
$ws.Range("N15").Value = 1.075614917723932
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.069177764164984
$ws.Range("D16").Value = 1.069204822376452
$ws.Range("E16").Value = 1.073070577514786
$ws.Range("F16").Value = 1.077619114296635
$ws.Range("I16").Value = 1.04711603007249
$ws.Range("J16").Value = 1.075274479559898
$ws.Range("K16").Value = 1.072528032482487
$ws.Range("L16").Value = 1.076380984965014
$ws.Range("M16").Value = 1.080914586798667
$ws.Range("N16").Value = 1.076801492656579
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.070054174208224
$ws.Range("D17").Value = 1.06990070137929
$ws.Range("E17").Value = 1.07392799558958
$ws.Range("F17").Value = 1.0783693421426
$ws.Range("I17").Value = 1.047363564015588
$ws.Range("J17").Value = 1.0760155098099
$ws.Range("K17").Value = 1.073152503629083
$ws.Range("L17").Value = 1.077166774766748
$ws.Range("M17").Value = 1.08159388276783
$ws.Range("N17").Value = 1.077543575254533
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.070564581762858
$ws.Range("D18").Value = 1.070305877193239
$ws.Range("E18").Value = 1.074427375047307
$ws.Range("F18").Value = 1.078806229970521
$ws.Range("I18").Value = 1.04750743181839
$ws.Range("J18").Value = 1.076446950683128
$ws.Range("K18").Value = 1.073515979517246
$ws.Range("L18").Value = 1.077624337060316
$ws.Range("M18").Value = 1.081989351433284
$ws.Range("N18").Value = 1.077975628823289
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.070738485162502
$ws.Range("D19").Value = 1.070443910709187
$ws.Range("E19").Value = 1.074597526501691
$ws.Range("F19").Value = 1.078955078425306
$ws.Range("I19").Value = 1.047556400129556
$ws.Range("J19").Value = 1.076593927982425
$ws.Range("K19").Value = 1.073639786088327
$ws.Range("L19").Value = 1.077780223497387
$ws.Range("M19").Value = 1.082124069203337
$ws.Range("N19").Value = 1.07812281484721
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.069960225457859
$ws.Range("D20").Value = 1.069826114701918
$ws.Range("E20").Value = 1.073836079352883
$ws.Range("F20").Value = 1.078288923197384
$ws.Range("I20").Value = 1.047337059281886
$ws.Range("J20").Value = 1.075936086189017
$ws.Range("K20").Value = 1.073085583388541
$ws.Range("L20").Value = 1.077082547439506
$ws.Range("M20").Value = 1.081521078866327
$ws.Range("N20").Value = 1.077464038842999
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.067420680824835
$ws.Range("D21").Value = 1.067809084305809
$ws.Range("E21").Value = 1.071351781516821
$ws.Range("F21").Value = 1.076114800601782
$ws.Range("I21").Value = 1.046617884258551
$ws.Range("J21").Value = 1.073788022078897
$ws.Range("K21").Value = 1.071274729494384
$ws.Range("L21").Value = 1.074805138607379
$ws.Range("M21").Value = 1.079551777083481
$ws.Range("N21").Value = 1.075312924235456
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.065816595804932
$ws.Range("D22").Value = 1.066534209610222
$ws.Range("E22").Value = 1.069782893372665
$ws.Range("F22").Value = 1.074741245034795
$ws.Range("I22").Value = 1.046160985753408
$ws.Range("J22").Value = 1.072430098940402
$ws.Range("K22").Value = 1.070129056414498
$ws.Range("L22").Value = 1.073366006341822
$ws.Range("M22").Value = 1.078306599948717
$ws.Range("N22").Value = 1.073953072690349
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.066667664021695
$ws.Range("D23").Value = 1.067210689171273
$ws.Range("E23").Value = 1.070615258664357
$ws.Range("F23").Value = 1.075470029329552
$ws.Range("I23").Value = 1.046403648493787
$ws.Range("J23").Value = 1.073150667991507
$ws.Range("K23").Value = 1.070737083945774
$ws.Range("L23").Value = 1.074129615496394
$ws.Range("M23").Value = 1.078967365784332
$ws.Range("N23").Value = 1.074674665032158
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.070002679286256
$ws.Range("D24").Value = 1.069859819430297
$ws.Range("E24").Value = 1.073877614622704
$ws.Range("F24").Value = 1.078325263241723
$ws.Range("I24").Value = 1.047349037219661
$ws.Range("J24").Value = 1.075971976742866
$ws.Range("K24").Value = 1.073115824135345
$ws.Range("L24").Value = 1.077120608538384
$ws.Range("M24").Value = 1.08155397814176
$ws.Range("N24").Value = 1.077499980365552
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.073842019746309
$ws.Range("D25").Value = 1.072905902257425
$ws.Range("E25").Value = 1.077634576580955
$ws.Range("F25").Value = 1.081610965348059
$ws.Range("I25").Value = 1.04842596691802
$ws.Range("J25").Value = 1.079215093560477
$ws.Range("K25").Value = 1.075846209370278
$ws.Range("L25").Value = 1.080561186202308
$ws.Range("M25").Value = 1.084526157316374
$ws.Range("N25").Value = 1.080747702781035
